$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-14.151289456726818, 2.1586155921367434, 1),
    @(-5.8697131578048243, 3.6195920279994089, 1),
    @(-2.1586155921367434, 5.8697131578048243, 1),
    @(-3.6195920279994089, 14.151289456726818, 1)
)

$row = 14
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
